$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$arr = New-Object 'double[,]' 20,10
$arr[0,0] = -18.1813755249159
$arr[0,1] = -0.7308005628491562
$arr[0,2] = -18.1813755249159
$arr[0,3] = -18.1813755249159
$arr[0,4] = -18.1813755249159
$arr[0,5] = -18.1813755249159
$arr[0,6] = -18.1813755249159
$arr[0,7] = -18.1813755249159
$arr[0,8] = -18.1813755249159
$arr[0,9] = -18.1813755249159
$arr[1,0] = -18.1813755249159
$arr[1,1] = -18.1813755249159
$arr[1,2] = -18.1813755249159
$arr[1,3] = -18.1813755249159
$arr[1,4] = -18.1813755249159
$arr[1,5] = -18.1813755249159
$arr[1,6] = -18.1813755249159
$arr[1,7] = -0.4552305561090912
$arr[1,8] = -18.1813755249159
$arr[1,9] = -18.1813755249159
$arr[2,0] = -18.1813755249159
$arr[2,1] = -0.5769308842683269
$arr[2,2] = -0.1601949122714088
$arr[2,3] = -18.1813755249159
$arr[2,4] = 4.022927250135273
$arr[2,5] = -18.1813755249159
$arr[2,6] = 2.96281280667059
$arr[2,7] = -18.1813755249159
$arr[2,8] = 3.487681232325651
$arr[2,9] = -18.1813755249159
$arr[3,0] = -18.1813755249159
$arr[3,1] = -0.07062206132905523
$arr[3,2] = -18.1813755249159
$arr[3,3] = -18.1813755249159
$arr[3,4] = -18.1813755249159
$arr[3,5] = 3.838864188237006
$arr[3,6] = -18.1813755249159
$arr[3,7] = -18.1813755249159
$arr[3,8] = -18.1813755249159
$arr[3,9] = -18.1813755249159
$arr[4,0] = -18.1813755249159
$arr[4,1] = -18.1813755249159
$arr[4,2] = -18.1813755249159
$arr[4,3] = -18.1813755249159
$arr[4,4] = -18.1813755249159
$arr[4,5] = -18.1813755249159
$arr[4,6] = -18.1813755249159
$arr[4,7] = -18.1813755249159
$arr[4,8] = -18.1813755249159
$arr[4,9] = -18.1813755249159
$arr[5,0] = 3.375735211870081
$arr[5,1] = -18.1813755249159
$arr[5,2] = -18.1813755249159
$arr[5,3] = -18.1813755249159
$arr[5,4] = -18.1813755249159
$arr[5,5] = -18.1813755249159
$arr[5,6] = -18.1813755249159
$arr[5,7] = -18.1813755249159
$arr[5,8] = -18.1813755249159
$arr[5,9] = -18.1813755249159
$arr[6,0] = -18.1813755249159
$arr[6,1] = -18.1813755249159
$arr[6,2] = -18.1813755249159
$arr[6,3] = 1.121765288477397
$arr[6,4] = -18.1813755249159
$arr[6,5] = -18.1813755249159
$arr[6,6] = -18.1813755249159
$arr[6,7] = -18.1813755249159
$arr[6,8] = -18.1813755249159
$arr[6,9] = -18.1813755249159
$arr[7,0] = 3.266027086506802
$arr[7,1] = -18.1813755249159
$arr[7,2] = -18.1813755249159
$arr[7,3] = -18.1813755249159
$arr[7,4] = -18.1813755249159
$arr[7,5] = -18.1813755249159
$arr[7,6] = -18.1813755249159
$arr[7,7] = -18.1813755249159
$arr[7,8] = -18.1813755249159
$arr[7,9] = -18.1813755249159
$arr[8,0] = -18.1813755249159
$arr[8,1] = -18.1813755249159
$arr[8,2] = -18.1813755249159
$arr[8,3] = -18.1813755249159
$arr[8,4] = -18.1813755249159
$arr[8,5] = -18.1813755249159
$arr[8,6] = -18.1813755249159
$arr[8,7] = -0.2443868457174131
$arr[8,8] = -18.1813755249159
$arr[8,9] = -18.1813755249159
$arr[9,0] = -18.1813755249159
$arr[9,1] = -18.1813755249159
$arr[9,2] = -18.1813755249159
$arr[9,3] = 2.628126500727879
$arr[9,4] = -18.1813755249159
$arr[9,5] = 1.227478693154117
$arr[9,6] = -18.1813755249159
$arr[9,7] = -18.1813755249159
$arr[9,8] = -18.1813755249159
$arr[9,9] = -18.1813755249159
$arr[10,0] = -18.1813755249159
$arr[10,1] = -18.1813755249159
$arr[10,2] = -18.1813755249159
$arr[10,3] = -18.1813755249159
$arr[10,4] = -18.1813755249159
$arr[10,5] = -18.1813755249159
$arr[10,6] = -18.1813755249159
$arr[10,7] = -18.1813755249159
$arr[10,8] = -18.1813755249159
$arr[10,9] = -18.1813755249159
$arr[11,0] = -18.1813755249159
$arr[11,1] = -18.1813755249159
$arr[11,2] = -18.1813755249159
$arr[11,3] = 2.019203442795556
$arr[11,4] = -18.1813755249159
$arr[11,5] = -18.1813755249159
$arr[11,6] = -18.1813755249159
$arr[11,7] = -18.1813755249159
$arr[11,8] = 0.4785263228898293
$arr[11,9] = 4.321923484266163
$arr[12,0] = -18.1813755249159
$arr[12,1] = -18.1813755249159
$arr[12,2] = 0.6459950930534166
$arr[12,3] = -18.1813755249159
$arr[12,4] = -18.1813755249159
$arr[12,5] = -18.1813755249159
$arr[12,6] = -18.1813755249159
$arr[12,7] = -18.1813755249159
$arr[12,8] = -18.1813755249159
$arr[12,9] = -18.1813755249159
$arr[13,0] = -18.1813755249159
$arr[13,1] = -18.1813755249159
$arr[13,2] = -0.574760667030876
$arr[13,3] = -18.1813755249159
$arr[13,4] = -18.1813755249159
$arr[13,5] = -18.1813755249159
$arr[13,6] = -18.1813755249159
$arr[13,7] = -18.1813755249159
$arr[13,8] = -18.1813755249159
$arr[13,9] = -18.1813755249159
$arr[14,0] = -18.1813755249159
$arr[14,1] = -18.1813755249159
$arr[14,2] = -18.1813755249159
$arr[14,3] = -18.1813755249159
$arr[14,4] = -18.1813755249159
$arr[14,5] = -18.1813755249159
$arr[14,6] = -18.1813755249159
$arr[14,7] = -18.1813755249159
$arr[14,8] = 1.116358005068239
$arr[14,9] = -18.1813755249159
$arr[15,0] = -18.1813755249159
$arr[15,1] = -0.05904799109232499
$arr[15,2] = -0.5648195021388452
$arr[15,3] = -18.1813755249159
$arr[15,4] = -18.1813755249159
$arr[15,5] = -18.1813755249159
$arr[15,6] = 1.615685914095297
$arr[15,7] = -0.9407611871189863
$arr[15,8] = 0.967153191156333
$arr[15,9] = -18.1813755249159
$arr[16,0] = -18.1813755249159
$arr[16,1] = -18.1813755249159
$arr[16,2] = -18.1813755249159
$arr[16,3] = -18.1813755249159
$arr[16,4] = -18.1813755249159
$arr[16,5] = -18.1813755249159
$arr[16,6] = 1.519594119958756
$arr[16,7] = -1.265577822703746
$arr[16,8] = 1.707589145175037
$arr[16,9] = -18.1813755249159
$arr[17,0] = -18.1813755249159
$arr[17,1] = -18.1813755249159
$arr[17,2] = 2.837914529207062
$arr[17,3] = -18.1813755249159
$arr[17,4] = -18.1813755249159
$arr[17,5] = -18.1813755249159
$arr[17,6] = 1.912945762031749
$arr[17,7] = -0.1629090824901449
$arr[17,8] = -18.1813755249159
$arr[17,9] = -18.1813755249159
$arr[18,0] = -18.1813755249159
$arr[18,1] = 3.417059140642603
$arr[18,2] = 3.176791123808811
$arr[18,3] = -18.1813755249159
$arr[18,4] = 1.90444052307381
$arr[18,5] = -18.1813755249159
$arr[18,6] = 0.1663885629224279
$arr[18,7] = 4.052787299929604
$arr[18,8] = -18.1813755249159
$arr[18,9] = -18.1813755249159
$arr[19,0] = -18.1813755249159
$arr[19,1] = 2.616650839378285
$arr[19,2] = -18.1813755249159
$arr[19,3] = 2.923717724133317
$arr[19,4] = -18.1813755249159
$arr[19,5] = 1.743832628910085
$arr[19,6] = 0.4686699719358235
$arr[19,7] = -18.1813755249159
$arr[19,8] = -18.1813755249159
$arr[19,9] = -18.1813755249159
$ws.Range("B2:K21").Value = $arr
